$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "linear" row (row 2) with new retrained ridge metrics ---
$ws.Range("B2").Value2 = "ridge"
$ws.Range("C2").Value2 = 244.39
$ws.Range("D2").Value2 = 93763.62
$ws.Range("E2").Value2 = 0.72
$ws.Range("F2").Value2 = 180.53
$ws.Range("G2").Value2 = 60758.2
$ws.Range("H2").Value2 = 0.79
$ws.Range("I2").Value2 = 131.16
$ws.Range("J2").Value2 = 29017.98
$ws.Range("K2").Value2 = 0.92

# --- Insert two new rows before the old row 3 (baseline-rent), pushing it down to row 5 ---
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Match the bordered/centered/bold style used by the rest of column A
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 3: rf model (new MLFlow run) ---
$ws.Range("A3").Value2 = 1
$ws.Range("B3").Value2 = "rf"
$ws.Range("C3").Value2 = 192.81
$ws.Range("D3").Value2 = 83040.97
$ws.Range("E3").Value2 = 0.75
$ws.Range("F3").Value2 = 138.6
$ws.Range("G3").Value2 = 32505.34
$ws.Range("H3").Value2 = 0.89
$ws.Range("I3").Value2 = 65.56999999999999
$ws.Range("J3").Value2 = 7753.23
$ws.Range("K3").Value2 = 0.98

# --- Row 4: lasso model (new MLFlow run) ---
$ws.Range("A4").Value2 = 2
$ws.Range("B4").Value2 = "lasso"
$ws.Range("C4").Value2 = 249.11
$ws.Range("D4").Value2 = 97304.32000000001
$ws.Range("E4").Value2 = 0.71
$ws.Range("F4").Value2 = 190.04
$ws.Range("G4").Value2 = 66804.00999999999
$ws.Range("H4").Value2 = 0.77
$ws.Range("I4").Value2 = 132.14
$ws.Range("J4").Value2 = 29166.11
$ws.Range("K4").Value2 = 0.92

# --- Row 5: baseline-rent, now A=3 (shifted down by the inserted rows) ---
$ws.Range("A5").Value2 = 3
